# "putting boosts around maps"
# Add boost items to the Pizzabox/Boost table (column C = extra boost),
# and upgrade the L14 cheese boost to "cheese armor".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# L34 row (row 16): add third-column boost "basil"
$ws.Range("C16").Value = "basil"

# L14 row (row 8): Boost becomes "cheese armor", new third-column boost "bacon"
$ws.Range("B8").Value = "cheese armor"
$ws.Range("C8").Value = "bacon"

# L24 row (row 12): add third-column boost "garlic bread"
$ws.Range("C12").Value = "garlic bread"

# L44 row (row 20): add third-column boost "cheese"
$ws.Range("C20").Value = "cheese"

# Widen column C to fit the new boost names, and nudge column B width
$ws.Columns.Item(2).ColumnWidth = 13.43
$ws.Columns.Item(3).ColumnWidth = 12.7

# Update the active selection to match the authored workbook
$ws.Range("F12").Select()
